# Applies the "Updated ESP model - 2025-08-18 17:18" edit to Sets-vervestacks.xlsx
# Target workbook has two sheets:
#   VEDA_Sets-Comm  (sheet1.xml)
#   VEDA_Sets-Proc  (sheet2.xml, active tab)

$wb = $excel.ActiveWorkbook

$wsComm = $wb.Worksheets.Item("VEDA_Sets-Comm")
$wsProc = $wb.Worksheets.Item("VEDA_Sets-Proc")

# ---------------------------------------------------------------------------
# VEDA_Sets-Proc (sheet2): update pattern strings + add T_Pos_AndOr /
# T_Neg_AndOr ("And"/"Or") cells on the rows that now carry two match
# patterns, add a new exclusion row for nuclear SMR, and add a brand new
# "Grid" set row.
# ---------------------------------------------------------------------------

# Row 3 (CCGT): widen pattern, duplicate SetName into SetDesc, add And/Or
$wsProc.Range("B3").Value = "ep_gas_combined_cycle*,ep_oil_combined_cycle*,CCGT*,*GasCC*"
$wsProc.Range("G3").Value = "CCGT"
$wsProc.Range("H3").Value = "And"
$wsProc.Range("I3").Value = "Or"

# Row 7 (OCGT (Peaker)): widen pattern, add And/Or
$wsProc.Range("B7").Value = "ep_gas_gas_turbine*,ep_oil_gas_turbine*,gas turbine*,EN*CT*"
$wsProc.Range("H7").Value = "And"
$wsProc.Range("I7").Value = "Or"

# Row 17 (Nuclear / ELE): add SMR exclusion pattern + And/Or
$wsProc.Range("B17").Value = "-*SMR"
$wsProc.Range("H17").Value = "And"
$wsProc.Range("I17").Value = "Or"

# Row 19 (Util Batt Stg): replace pattern, add And/Or
$wsProc.Range("B19").Value = "EN*STG?hb*,-*EV*"
$wsProc.Range("H19").Value = "And"
$wsProc.Range("I19").Value = "Or"

# Row 21 (new): Grid / IRE set
$wsProc.Range("A21").Value = "IRE"
$wsProc.Range("B21").Value = "g[_]*"
$wsProc.Range("F21").Value = "Grid"

# Final selection on VEDA_Sets-Proc lands on the newly added B21
$wsProc.Range("B21").Select()

# ---------------------------------------------------------------------------
# VEDA_Sets-Comm (sheet1): reset the lingering A3:H9 selection back to A1.
# ---------------------------------------------------------------------------
$wsComm.Range("A1").Select()

# Restore VEDA_Sets-Proc as the active sheet/tab
$wsProc.Activate()
